# Update SNMP page and fresh write Certificate
#
# The "SNMP" worksheet (3rd sheet, sheet3.xml) gets 6 new columns (C..H)
# added to its 2-row table, a couple of existing cells rewritten, column
# widths / bestFit sizing applied, and the active selection moved to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNMP")
$ws.Activate() | Out-Null

# --- Row 1 (headers) -------------------------------------------------
# A1 "securityName" and B1 "authPro" already exist; add the rest.
$ws.Cells.Item(1, 3).Value = "authPass"
$ws.Cells.Item(1, 4).Value = "privPro"
$ws.Cells.Item(1, 5).Value = "privPass"

# --- Row 2 (values) ---------------------------------------------------
# A2 "khang" stays as-is.
$ws.Cells.Item(2, 4).Value = "DES|AES128|AES192|AES256"
$ws.Cells.Item(2, 2).Value = "MD5|SHA"

# --- Row 1 (more headers) ---------------------------------------------
$ws.Cells.Item(1, 6).Value = "SNMPserver"
$ws.Cells.Item(1, 7).Value = "SNMPuser"
$ws.Cells.Item(1, 8).Value = "SNMPpass"

# --- Row 2 (more values) ----------------------------------------------
$ws.Cells.Item(2, 7).Value = "root"
$ws.Cells.Item(2, 6).Value = "100.30.7.130"
$ws.Cells.Item(2, 8).Value = "1_Abc_123"
$ws.Cells.Item(2, 3).Value = "1_Abc_123"
$ws.Cells.Item(2, 5).Value = "1_Abc_123"

# --- Column widths (best-fit look) ------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.428571428571429
$ws.Columns.Item(2).ColumnWidth = 7.857142857142857
$ws.Columns.Item(3).ColumnWidth = 9.285714285714286
$ws.Columns.Item(4).ColumnWidth = 25.571428571428573
$ws.Columns.Item(5).ColumnWidth = 9.285714285714286
$ws.Columns.Item(6).ColumnWidth = 10.714285714285714
$ws.Columns.Item(7).ColumnWidth = 8.857142857142858
$ws.Columns.Item(8).ColumnWidth = 9.285714285714286

# --- Selection ----------------------------------------------------------
$ws.Range("F2").Select() | Out-Null
